$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 3 = "logic for computer science"
# Clear the last-updated date in B3 but keep its (date) number format
$ws.Range("B3").ClearContents()

# Mark Lec (column C) as Done
$ws.Range("C3").Value = "Done"

# Remove the old Rec/HW zero markers (column D/E) for this row
$ws.Range("D3").ClearContents()
$ws.Range("E3").ClearContents()

# Add a new note in column F
$ws.Range("F3").Value = "Complete gedels proof"

# Move the active selection from G4 to F4
$ws.Range("F4").Select()
